$wb = $excel.ActiveWorkbook

# --- Sheet "KPI nhan vien": insert 4 new detail rows (14-17) before the END row ---
$ws1 = $wb.Worksheets.Item("KPI nhan vien")

[void]$ws1.Rows("14:17").Insert()
$ws1.Rows("14:17").OutlineLevel = 1

$ws1.Range("A14:A17").Value = "{Mã NV}"
$ws1.Range("B14:B17").Value = "{Tên nhân viên}"

$ws1.Range("C14").Value = "SKU/ Đơn hàng trực tiếp"
$ws1.Range("C15").Value = "Doanh thu đơn hàng trực tiếp"
$ws1.Range("C16").Value = "Tổng sản lượng đơn hàng trực tiếp"
$ws1.Range("C17").Value = "Số đơn hàng trực tiếp"

# --- Sheet "Chi tieu": fill the existing blank cells + extend with 4 new rows ---
$ws2 = $wb.Worksheets.Item("Chi tieu")

$ws2.Range("A2").Value = "Số lần viếng thăm đại lý"
$ws2.Range("A3").Value = "Số đại lý tạo mới"
$ws2.Range("A4").Value = "Số đại lý viếng thăm"
$ws2.Range("A5").Value = "SKU/ Đơn hàng gián tiếp"
$ws2.Range("A6").Value = "Doanh thu đơn hàng gián tiếp"
$ws2.Range("A7").Value = "Tổng sản lượng đơn hàng gián tiếp"
$ws2.Range("A8").Value = "Số đơn hàng gián tiếp"

$ws2.Range("A9").Value = "SKU/ Đơn hàng trực tiếp"
$ws2.Range("A10").Value = "Doanh thu đơn hàng trực tiếp"
$ws2.Range("A11").Value = "Tổng sản lượng đơn hàng trực tiếp"
$ws2.Range("A12").Value = "Số đơn hàng trực tiếp"

# Copy the existing row-8 format down onto the 4 freshly appended rows so they
# pick up the same cell style (s="1") as the rest of the "Chi tieu" list.
[void]$ws2.Range("A8").Copy()
[void]$ws2.Range("A9:A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the recorded selections in the two touched sheets.
[void]$ws2.Range("A2:A12").Select()

[void]$ws1.Activate()
[void]$ws1.Range("F21").Select()
